# Auto update Excel log
# Adds new sensor log rows to ALERTS, Proximity, mmWave, and Camera sheets.

$wb = $excel.ActiveWorkbook

function Add-LogRow($ws, $row, $values) {
    $lastCol = [char]([int][char]'A' + $values.Length - 1)
    $rng = $ws.Range("A" + $row + ":" + $lastCol + $row)

    # Force text number format so date/time-looking strings (e.g. "2026-02-01",
    # "14:40:33") are stored as plain text instead of being auto-converted to
    # date/time serial numbers.
    $rng.NumberFormat = "@"

    $arr = New-Object 'object[,]' 1,$values.Length
    for ($i = 0; $i -lt $values.Length; $i++) {
        $arr[0,$i] = $values[$i]
    }
    $rng.Value = $arr

    # Reset style back to Normal so no stray text-format style is left behind
    # on the cells (keeps values as text while matching default styling).
    $rng.Style = "Normal"
}

# --- ALERTS sheet: add rows 11-12 ---
$wsAlerts = $wb.Worksheets.Item("ALERTS")
Add-LogRow $wsAlerts 11 @("2026-02-01", "14:40:33", "14:00", "Living Room", "CRITICAL", "FALL_DETECTED")
Add-LogRow $wsAlerts 12 @("2026-02-01", "14:40:34", "14:00", "Living Room", "CRITICAL", "FALL_DETECTED")

# --- Proximity sheet: add row 31 ---
$wsProximity = $wb.Worksheets.Item("Proximity")
Add-LogRow $wsProximity 31 @("2026-02-01", "14:40:53", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")

# --- mmWave sheet: add rows 11-13 ---
$wsMmWave = $wb.Worksheets.Item("mmWave")
Add-LogRow $wsMmWave 11 @("2026-02-01", "14:40:12", "14:00", "Living Room", "NO_MOTION_DETECTED", "Inactive")
Add-LogRow $wsMmWave 12 @("2026-02-01", "14:40:42", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")
Add-LogRow $wsMmWave 13 @("2026-02-01", "14:40:52", "14:00", "Living Room", "PRESENCE_DETECTED", "Active")

# --- Camera sheet: add rows 18-19 ---
$wsCamera = $wb.Worksheets.Item("Camera")
Add-LogRow $wsCamera 18 @("2026-02-01", "14:40:55", "14:00", "Living Room Main Door", "Image Captured", "Active")
Add-LogRow $wsCamera 19 @("2026-02-01", "14:40:55", "14:00", "Living Room Main Door", "Image Received", "Active")

Write-Host "Rows added successfully."
